$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values keyed by cell address, derived from the corrected Q calculation
$newValues = @{
    "B2" = -34.07235526768176
    "I2" = -39.27027815818423
    "J2" = -0.09180947066840872
    "L2" = 42.70727241941287
    "M2" = 42.70727241941287
    "N2" = 116.4743793256715
    "O2" = 42.70727241941287
    "B3" = -16.84468718035663
    "I3" = -39.76939137462676
    "J3" = -0.09106907171138801
    "L3" = 39.78120303934642
    "M3" = 39.78120303934642
    "N3" = 108.4941901073084
    "O3" = 39.78120303934642
    "B4" = -4.00261595342505
    "I4" = -39.21276070314079
    "J4" = -0.0858862790123851
    "L4" = 37.43217558786801
    "M4" = 37.43217558786801
    "N4" = 102.0877516032764
    "O4" = 37.43217558786802
    "B5" = 1.758325662124207
    "I5" = -37.75910917540352
    "J5" = -0.08681177770861837
    "L5" = 36.18100620145194
    "M5" = 36.18100620145194
    "N5" = 98.6754714585053
    "O5" = 36.18100620145194
    "B6" = 3.87465396766234
    "I6" = -36.68639292369173
    "J6" = -0.0866266779693774
    "L6" = 35.63411433530171
    "M6" = 35.63411433530171
    "N6" = 97.18394818718646
    "O6" = 35.63411433530172
    "B7" = 3.87465396766234
    "I7" = -36.68639292369173
    "J7" = -0.0866266779693774
    "L7" = 35.63411433530171
    "M7" = 35.63411433530171
    "N7" = 97.18394818718646
    "O7" = 35.63411433530172
    "B8" = 2.92670076625933
    "I8" = -36.48649619468827
    "J8" = -0.08366508214135138
    "L8" = 35.76949859250323
    "M8" = 35.76949859250323
    "N8" = 97.55317797955425
    "O8" = 35.76949859250323
    "B9" = -8.3802791329434
    "I9" = -36.17910782839053
    "J9" = -0.08403528161986173
    "L9" = 37.71702965267636
    "M9" = 37.71702965267636
    "N9" = 102.864626325481
    "O9" = 37.71702965267636
    "B10" = -34.76599320240507
    "I10" = -37.18707765529936
    "J10" = -0.0858862790123851
    "L10" = 42.49214764958622
    "M10" = 42.49214764958622
    "N10" = 115.8876754079624
    "O10" = 42.49214764958623
    "B11" = -60.29752499834723
    "I11" = -37.89080214624073
    "J11" = -0.08088858605259475
    "L11" = 47.00639194022902
    "M11" = 47.00639194022902
    "N11" = 128.1992507460791
    "O11" = 47.00639194022904
    "B12" = -72.06528231819937
    "I12" = -37.58184226477783
    "J12" = -1.057104610861359
    "L12" = 49.13563824158147
    "M12" = 49.13563824158147
    "N12" = 134.006286113404
    "O12" = 49.13563824158146
    "B13" = -75.35929950159148
    "I13" = -37.23736614058328
    "J13" = -7.595197600678148
    "L13" = 50.78015535377907
    "M13" = 50.78015535377907
    "N13" = 138.4913327830338
    "O13" = 50.78015535377909
    "B14" = -71.5462281690246
    "I14" = -36.74171025750397
    "J14" = -13.42287779125041
    "L14" = 51.06745401824455
    "M14" = 51.06745401824455
    "N14" = 139.2748745952124
    "O14" = 51.06745401824452
    "B15" = -67.43637480566525
    "I15" = -36.00561255416113
    "J15" = -15.74495402015174
    "L15" = 50.65326899062375
    "M15" = 50.65326899062375
    "N15" = 138.1452790653375
    "O15" = 50.65326899062376
    "B16" = -65.90082688938367
    "I16" = -36.72505219675372
    "J16" = -16.60696350584271
    "L16" = 50.65722691387813
    "M16" = 50.65722691387813
    "N16" = 138.1560734014858
    "O16" = 50.65722691387813
    "B17" = -64.51266693973224
    "I17" = -37.65916081090182
    "J17" = -17.25
    "L17" = 50.68248777700171
    "M17" = 50.68248777700171
    "N17" = 138.2249666645501
    "O17" = 50.68248777700174
    "B18" = -57.31750392205134
    "I18" = -37.56895584042388
    "J18" = -16.67711630701876
    "L18" = 49.33714309903274
    "M18" = 49.33714309903274
    "N18" = 134.5558448155438
    "O18" = 49.33714309903272
    "B19" = -54.12644504593231
    "I19" = -37.98697887434616
    "J19" = -14.63768737995343
    "L19" = 48.50167206383475
    "M19" = 48.50167206383475
    "N19" = 132.277287446822
    "O19" = 48.50167206383475
    "B20" = -61.764826859816
    "I20" = -37.9049457827268
    "J20" = -7.235178607835337
    "L20" = 48.50784176773129
    "M20" = 48.50784176773129
    "N20" = 132.2941139119944
    "O20" = 48.50784176773129
    "B21" = -80.84133812456298
    "I21" = -38.37765754439528
    "J21" = -0.6650633631281266
    "L21" = 50.68493237665884
    "M21" = 50.68493237665884
    "N21" = 138.2316337545241
    "O21" = 50.68493237665882
    "B22" = -92.88025623809939
    "I22" = -39.44408773544666
    "J22" = -0.1155022372925032
    "L22" = 52.79264291912713
    "M22" = 52.79264291912713
    "N22" = 143.9799352339831
    "O22" = 52.79264291912716
    "B23" = -101.375560942861
    "I23" = -40.59443683631548
    "J23" = -0.1084684472009485
    "L23" = 54.40235358858764
    "M23" = 54.40235358858764
    "N23" = 148.3700552416026
    "O23" = 54.40235358858763
    "B24" = -104.176868333314
    "I24" = -41.40376714635295
    "J24" = -0.1012495573701813
    "L24" = 55
    "M24" = 55
    "N24" = 150
    "O24" = 55
    "B25" = -84.97308619446042
    "I25" = -42
    "J25" = -0.09884326075990657
    "L25" = 51.86753656316803
    "M25" = 51.86753656316803
    "N25" = 141.4569178995492
    "O25" = 51.86753656316802
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
